$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 4
$ws.Range("H4").Value = 550
$ws.Range("I4").Value = 400
$ws.Range("J4").Value = 580
$ws.Range("K4").Value = 400
$ws.Range("L4").Value = 580
$ws.Range("M4").Value = -286
$ws.Range("N4").Value = -808

# Row 13
$ws.Range("H13").Value = 35000
$ws.Range("I13").Value = 50000
$ws.Range("J13").Value = 20000
$ws.Range("K13").Value = 50000
$ws.Range("L13").Value = 20000
$ws.Range("M13").Value = -49831
$ws.Range("N13").Value = -20338

# Row 28
$ws.Range("H28").Value = 3740.6667
$ws.Range("I28").Value = 3414.6667
$ws.Range("J28").Value = 4066.6667
$ws.Range("K28").Value = 3414.6667
$ws.Range("L28").Value = 4066.6667
$ws.Range("M28").Value = -2929.6667
$ws.Range("N28").Value = -5036.6667

# Row 88
$ws.Range("H88").Value = 2153.0667
$ws.Range("J88").Value = 2021.1428
$ws.Range("L88").Value = 2021.1428
$ws.Range("N88").Value = -2833.1428

# Row 91
$ws.Range("H91").Value = 2153.0667
$ws.Range("J91").Value = 2021.1428
$ws.Range("L91").Value = 2021.1428
$ws.Range("N91").Value = -4829.1428

# Row 98
$ws.Range("H98").Value = 9576.333000000001
$ws.Range("J98").Value = 25253
$ws.Range("L98").Value = 25253
$ws.Range("N98").Value = -28249

# Row 112
$ws.Range("H112").Value = 1927.7742
$ws.Range("J112").Value = 2160.5
$ws.Range("L112").Value = 6481.5
$ws.Range("N112").Value = -8697.5

# Row 122
$ws.Range("H122").Value = 9576.333000000001
$ws.Range("J122").Value = 25253
$ws.Range("L122").Value = 75759
$ws.Range("N122").Value = -80659

# Row 127
$ws.Range("H127").Value = 1504.7333
$ws.Range("I127").Value = 721.5
$ws.Range("J127").Value = 2399.8572
$ws.Range("K127").Value = 2164.5
$ws.Range("L127").Value = 7199.571599999999
$ws.Range("M127").Value = 2795.5
$ws.Range("N127").Value = -17119.5716

# Row 132
$ws.Range("H132").Value = 1878.1356
$ws.Range("I132").Value = 1536.9767
$ws.Range("J132").Value = 2795
$ws.Range("K132").Value = 4610.9301
$ws.Range("L132").Value = 8385
$ws.Range("M132").Value = -2080.9301
$ws.Range("N132").Value = -13445

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 12400.918
$ws.Range("I32").Value = 13799.288
$ws.Range("K32").Value = 13799.288
$ws.Range("M32").Value = -13512.288

# Row 74
$ws.Range("H74").Value = 1132.0834
$ws.Range("I74").Value = 713.73334
$ws.Range("K74").Value = 713.73334
$ws.Range("M74").Value = 160.26666

# Row 77
$ws.Range("H77").Value = 1132.0834
$ws.Range("I77").Value = 713.73334
$ws.Range("K77").Value = 3568.6667
$ws.Range("M77").Value = 799.3333000000002

# Row 88
$ws.Range("H88").Value = 2442.3333
$ws.Range("J88").Value = 2570.7144
$ws.Range("L88").Value = 2570.7144
$ws.Range("N88").Value = -3382.7144

# Row 91
$ws.Range("H91").Value = 2442.3333
$ws.Range("J91").Value = 2570.7144
$ws.Range("L91").Value = 2570.7144
$ws.Range("N91").Value = -5378.7144

# Row 123
$ws.Range("H123").Value = 30429
$ws.Range("J123").Value = 30429
$ws.Range("L123").Value = 30429
$ws.Range("N123").Value = -40229

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 251449.75
$ws.Range("I20").Value = 251449.75
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 251449.75
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -251202.75
$ws.Range("N20").ClearContents()

# Row 26
$ws.Range("H26").Value = 258147.5
$ws.Range("I26").Value = 258147.5
$ws.Range("K26").Value = 258147.5
$ws.Range("M26").Value = -257855.5

# Row 68
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

# Row 71
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

# Row 96
$ws.Range("H96").Value = 89788.84
$ws.Range("I96").Value = 107584.2
$ws.Range("K96").Value = 107584.2
$ws.Range("M96").Value = -104838.2

# Row 105
$ws.Range("H105").Value = 7145856.5
$ws.Range("I105").Value = 12990003
$ws.Range("J105").Value = 3011.6667
$ws.Range("K105").Value = 12990003
$ws.Range("L105").Value = 3011.6667
$ws.Range("M105").Value = -12988256
$ws.Range("N105").Value = -6505.6667

$ws = $wb.Worksheets.Item("CUL")
# Row 33
$ws.Range("H33").Value = 936.7143
$ws.Range("I33").Value = 575.25
$ws.Range("J33").Value = 1081.3
$ws.Range("K33").Value = 3451.5
$ws.Range("L33").Value = 6487.799999999999
$ws.Range("M33").Value = -3168.5
$ws.Range("N33").Value = -7053.799999999999

# Row 107
$ws.Range("H107").Value = 2000
$ws.Range("J107").Value = 2000
$ws.Range("L107").Value = 6000
$ws.Range("N107").Value = -9840

# Row 121
$ws.Range("H121").Value = 6100.727
$ws.Range("J121").Value = 7812.25
$ws.Range("L121").Value = 23436.75
$ws.Range("N121").Value = -26056.75

# Row 131
$ws.Range("H131").Value = 19958.562
$ws.Range("I131").Value = 338.17648
$ws.Range("J131").Value = 30718.129
$ws.Range("K131").Value = 1014.52944
$ws.Range("L131").Value = 92154.387
$ws.Range("M131").Value = 4025.47056
$ws.Range("N131").Value = -102234.387

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 2744.2666
$ws.Range("I102").Value = 2586.3333
$ws.Range("J102").Value = 3376
$ws.Range("K102").Value = 2586.3333
$ws.Range("L102").Value = 3376
$ws.Range("M102").Value = -964.3332999999998
$ws.Range("N102").Value = -6620

# Row 109
$ws.Range("H109").Value = 9068.929
$ws.Range("J109").Value = 9068.929
$ws.Range("L109").Value = 9068.929
$ws.Range("N109").Value = -11148.929

# Row 132
$ws.Range("H132").Value = 2394.4736
$ws.Range("I132").Value = 1874.625
$ws.Range("J132").Value = 2772.5454
$ws.Range("K132").Value = 5623.875
$ws.Range("L132").Value = 8317.636200000001
$ws.Range("M132").Value = -3093.875
$ws.Range("N132").Value = -13377.6362

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 6011.6665
$ws.Range("I7").Value = 4183.3335
$ws.Range("J7").Value = 9668.333000000001
$ws.Range("K7").Value = 4183.3335
$ws.Range("L7").Value = 9668.333000000001
$ws.Range("M7").Value = -4071.3335
$ws.Range("N7").Value = -9892.333000000001

# Row 61
$ws.Range("H61").Value = 3025
$ws.Range("I61").Value = 3287.5
$ws.Range("J61").Value = 2500
$ws.Range("K61").Value = 3287.5
$ws.Range("L61").Value = 2500
$ws.Range("M61").Value = -3085.5
$ws.Range("N61").Value = -2904

# Row 113
$ws.Range("H113").Value = 3025
$ws.Range("I113").Value = 3287.5
$ws.Range("J113").Value = 2500
$ws.Range("K113").Value = 3287.5
$ws.Range("L113").Value = 2500
$ws.Range("M113").Value = -1117.5
$ws.Range("N113").Value = -6840

# Row 126
$ws.Range("H126").Value = 6011.6665
$ws.Range("I126").Value = 4183.3335
$ws.Range("J126").Value = 9668.333000000001
$ws.Range("K126").Value = 12550.0005
$ws.Range("L126").Value = 29004.999
$ws.Range("M126").Value = -10080.0005
$ws.Range("N126").Value = -33944.999

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 3713.3333
$ws.Range("I62").Value = 3106
$ws.Range("K62").Value = 3106
$ws.Range("M62").Value = -2482

# Row 65
$ws.Range("H65").Value = 3713.3333
$ws.Range("I65").Value = 3106
$ws.Range("K65").Value = 15530
$ws.Range("M65").Value = -12410

# Row 123
$ws.Range("H123").Value = 22427.95
$ws.Range("J123").Value = 22427.95
$ws.Range("L123").Value = 22427.95
$ws.Range("N123").Value = -32227.95
